$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D6").Value = "Finished"
$ws.Range("E2:E6").Value = 0
